$wb = $excel.ActiveWorkbook

# ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 958276.25
$ws.Range("I11").Value = 958276.25
$ws.Range("K11").Value = 958276.25
$ws.Range("M11").Value = -958136.25

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1433.3334
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1433.3334
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1433.3334
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1783.3334

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2231.85
$ws.Range("J112").Value = 2414
$ws.Range("L112").Value = 7242
$ws.Range("N112").Value = -9458

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 12750

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1007.6
$ws.Range("J129").Value = 1047.5714
$ws.Range("L129").Value = 3142.7142
$ws.Range("N129").Value = -13142.7142

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3215.9648
$ws.Range("I137").Value = 2774.756
$ws.Range("J137").Value = 4346.5625
$ws.Range("K137").Value = 8324.268
$ws.Range("L137").Value = 13039.6875
$ws.Range("M137").Value = -5774.268
$ws.Range("N137").Value = -18139.6875

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3824.0908
$ws.Range("I138").Value = 1702.9412
$ws.Range("J138").Value = 4425.0835
$ws.Range("K138").Value = 5108.8236
$ws.Range("L138").Value = 13275.2505
$ws.Range("M138").Value = 31.17640000000029
$ws.Range("N138").Value = -23555.2505

# ARM row 10
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 17776.5
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 17776.5
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 17776.5
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -18116.5

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1033.2667
$ws.Range("I45").Value = 949.9
$ws.Range("K45").Value = 949.9
$ws.Range("M45").Value = -572.9

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3891.5833
$ws.Range("I122").Value = 1949.8334
$ws.Range("K122").Value = 5849.5002
$ws.Range("M122").Value = -3399.5002

# BSM row 42
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 77300
$ws.Range("J42").Value = 77300
$ws.Range("L42").Value = 77300
$ws.Range("N42").Value = -77956

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3610.524
$ws.Range("I99").Value = 1649.8334
$ws.Range("J99").Value = 4394.8
$ws.Range("K99").Value = 1649.8334
$ws.Range("L99").Value = 4394.8
$ws.Range("M99").Value = -151.8334
$ws.Range("N99").Value = -7390.8

# BSM row 109
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 39000
$ws.Range("J109").Value = 39000
$ws.Range("L109").Value = 39000
$ws.Range("N109").Value = -41774

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1357.12
$ws.Range("I16").Value = 1149.3125
$ws.Range("K16").Value = 1149.3125
$ws.Range("M16").Value = -862.3125

# CRP row 57
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 50000
$ws.Range("J57").Value = 50000
$ws.Range("L57").Value = 50000
$ws.Range("N57").Value = -51120

# CRP row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 12876.2
$ws.Range("I69").Value = 12876.2
$ws.Range("K69").Value = 12876.2
$ws.Range("M69").Value = -12127.2

# CRP row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 12876.2
$ws.Range("I72").Value = 12876.2
$ws.Range("K72").Value = 38628.60000000001
$ws.Range("M72").Value = -34884.60000000001

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1357.12
$ws.Range("I113").Value = 1149.3125
$ws.Range("K113").Value = 1149.3125
$ws.Range("M113").Value = 1020.6875

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2660.5
$ws.Range("I132").Value = 2245.825
$ws.Range("J132").Value = 4319.2
$ws.Range("K132").Value = 6737.474999999999
$ws.Range("L132").Value = 12957.6
$ws.Range("M132").Value = -4207.474999999999
$ws.Range("N132").Value = -18017.6

# CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1744
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 1744
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 5232
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -6094

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 35717144
$ws.Range("I92").Value = 726.8333
$ws.Range("J92").Value = 62504460
$ws.Range("K92").Value = 2180.4999
$ws.Range("L92").Value = 187513380
$ws.Range("M92").Value = -932.4998999999998
$ws.Range("N92").Value = -187515876

# CUL row 96
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 4080
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 4600
$ws.Range("K96").Value = 6000
$ws.Range("L96").Value = 13800
$ws.Range("M96").Value = -3941
$ws.Range("N96").Value = -17918

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1191.6666
$ws.Range("I98").Value = 1272.7273
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 3818.1819
$ws.Range("L98").Value = 900
$ws.Range("M98").Value = -2320.1819
$ws.Range("N98").Value = -3896

# GSM row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# GSM row 26
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 6668
$ws.Range("J26").Value = 8694.666999999999
$ws.Range("L26").Value = 8694.666999999999
$ws.Range("N26").Value = -9254.666999999999

# GSM row 50
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 6668
$ws.Range("J50").Value = 8694.666999999999
$ws.Range("L50").Value = 8694.666999999999
$ws.Range("N50").Value = -9690.666999999999

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 35716644
$ws.Range("J80").Value = 2749.5
$ws.Range("L80").Value = 2749.5
$ws.Range("N80").Value = -4745.5

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 35716644
$ws.Range("J83").Value = 2749.5
$ws.Range("L83").Value = 13747.5
$ws.Range("N83").Value = -23731.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1512
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1520
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1520
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -5860

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4160.4185
$ws.Range("J126").Value = 5178.174
$ws.Range("L126").Value = 15534.522
$ws.Range("N126").Value = -20474.522

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1084.0834
$ws.Range("I61").Value = 875.5
$ws.Range("J61").Value = 1501.25
$ws.Range("K61").Value = 875.5
$ws.Range("L61").Value = 1501.25
$ws.Range("M61").Value = -673.5
$ws.Range("N61").Value = -1905.25

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1084.0834
$ws.Range("I113").Value = 875.5
$ws.Range("J113").Value = 1501.25
$ws.Range("K113").Value = 875.5
$ws.Range("L113").Value = 1501.25
$ws.Range("M113").Value = 1294.5
$ws.Range("N113").Value = -5841.25

# WVR row 7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 50054.5
$ws.Range("I7").Value = 104
$ws.Range("J7").Value = 100005
$ws.Range("K7").Value = 104
$ws.Range("L7").Value = 100005
$ws.Range("M7").Value = 9
$ws.Range("N7").Value = -100231

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 10358.8
$ws.Range("I113").Value = 14584.286
$ws.Range("J113").Value = 499.33334
$ws.Range("K113").Value = 43752.858
$ws.Range("L113").Value = 1498.00002
$ws.Range("M113").Value = -41582.858
$ws.Range("N113").Value = -5838.000019999999

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 29608.428
$ws.Range("J123").Value = 29608.428
$ws.Range("L123").Value = 29608.428
$ws.Range("N123").Value = -39408.428

Write-Output "applied all changes"